$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row mirrors the formatting of the row above it (A: style 2, B: style 1)
$ws.Range("A86:B86").Copy()
$ws.Range("A87:B87").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A87").Value = "MigrationInfoMessage"
$ws.Range("B87").Value = "Sadece sinif ve ogrenci bilgileri aktarilacaktir. Eski  sinav ve optik formlar yeni sistemle uyumlu olmadigi icin bunlar aktarilmayacaktir."
